$d = $word.ActiveDocument
$t = $d.Tables(1)

$t.Cell(1,1).Range.Text = "76+8="
$t.Cell(1,2).Range.Text = "74-17="
$t.Cell(1,3).Range.Text = "95-78="
$t.Cell(1,4).Range.Text = "35-18="
$t.Cell(1,5).Range.Text = "62+36="
$t.Cell(2,1).Range.Text = "7+86="
$t.Cell(2,2).Range.Text = "94-27="
$t.Cell(2,3).Range.Text = "55+32="
$t.Cell(2,4).Range.Text = "82-53="
$t.Cell(2,5).Range.Text = "68-1="
$t.Cell(3,1).Range.Text = "17+57="
$t.Cell(3,2).Range.Text = "26+45="
$t.Cell(3,3).Range.Text = "42+30="
$t.Cell(3,4).Range.Text = "35+11="
$t.Cell(3,5).Range.Text = "63-48="
$t.Cell(4,1).Range.Text = "14+35="
$t.Cell(4,2).Range.Text = "94-26="
$t.Cell(4,3).Range.Text = "63+10="
$t.Cell(4,4).Range.Text = "22+77="
$t.Cell(4,5).Range.Text = "59-57="
$t.Cell(5,1).Range.Text = "44+1="
$t.Cell(5,2).Range.Text = "24+46="
$t.Cell(5,3).Range.Text = "89-18="
$t.Cell(5,4).Range.Text = "81-60="
$t.Cell(5,5).Range.Text = "85+11="
$t.Cell(6,1).Range.Text = "85-59="
$t.Cell(6,2).Range.Text = "79-14="
$t.Cell(6,3).Range.Text = "91-11="
$t.Cell(6,4).Range.Text = "92-62="
$t.Cell(6,5).Range.Text = "33-26="
$t.Cell(7,1).Range.Text = "14+67="
$t.Cell(7,2).Range.Text = "38+46="
$t.Cell(7,3).Range.Text = "30+51="
$t.Cell(7,4).Range.Text = "69-39="
$t.Cell(7,5).Range.Text = "13+50="
$t.Cell(8,1).Range.Text = "13+85="
$t.Cell(8,2).Range.Text = "51+13="
$t.Cell(8,3).Range.Text = "82-24="
$t.Cell(8,4).Range.Text = "96-94="
$t.Cell(8,5).Range.Text = "73-73="
$t.Cell(9,1).Range.Text = "36-20="
$t.Cell(9,2).Range.Text = "99-61="
$t.Cell(9,3).Range.Text = "15-2="
$t.Cell(9,4).Range.Text = "11+86="
$t.Cell(9,5).Range.Text = "85-57="
$t.Cell(10,1).Range.Text = "78-27="
$t.Cell(10,2).Range.Text = "30+45="
$t.Cell(10,3).Range.Text = "95-34="
$t.Cell(10,4).Range.Text = "50-9="
$t.Cell(10,5).Range.Text = "37-0="
$t.Cell(11,1).Range.Text = "40-13="
$t.Cell(11,2).Range.Text = "57-42="
$t.Cell(11,3).Range.Text = "44+54="
$t.Cell(11,4).Range.Text = "60-7="
$t.Cell(11,5).Range.Text = "89+10="
$t.Cell(12,1).Range.Text = "16-7="
$t.Cell(12,2).Range.Text = "27+50="
$t.Cell(12,3).Range.Text = "87-62="
$t.Cell(12,4).Range.Text = "81-24="
$t.Cell(12,5).Range.Text = "39+58="
$t.Cell(13,1).Range.Text = "15+16="
$t.Cell(13,2).Range.Text = "14+27="
$t.Cell(13,3).Range.Text = "5+33="
$t.Cell(13,4).Range.Text = "58+12="
$t.Cell(13,5).Range.Text = "84-64="
$t.Cell(14,1).Range.Text = "89+3="
$t.Cell(14,2).Range.Text = "53-44="
$t.Cell(14,3).Range.Text = "27+50="
$t.Cell(14,4).Range.Text = "51-12="
$t.Cell(14,5).Range.Text = "33-27="
$t.Cell(15,1).Range.Text = "62-34="
$t.Cell(15,2).Range.Text = "56-37="
$t.Cell(15,3).Range.Text = "62+35="
$t.Cell(15,4).Range.Text = "21+67="
$t.Cell(15,5).Range.Text = "15-12="
$t.Cell(16,1).Range.Text = "80-20="
$t.Cell(16,2).Range.Text = "30+19="
$t.Cell(16,3).Range.Text = "79-63="
$t.Cell(16,4).Range.Text = "86-12="
$t.Cell(16,5).Range.Text = "25+21="
$t.Cell(17,1).Range.Text = "41+13="
$t.Cell(17,2).Range.Text = "48+41="
$t.Cell(17,3).Range.Text = "13+74="
$t.Cell(17,4).Range.Text = "73+9="
$t.Cell(17,5).Range.Text = "69-24="
$t.Cell(18,1).Range.Text = "41-33="
$t.Cell(18,2).Range.Text = "24-9="
$t.Cell(18,3).Range.Text = "14+70="
$t.Cell(18,4).Range.Text = "8+30="
$t.Cell(18,5).Range.Text = "33+2="
$t.Cell(19,1).Range.Text = "64-37="
$t.Cell(19,2).Range.Text = "78-43="
$t.Cell(19,3).Range.Text = "44+11="
$t.Cell(19,4).Range.Text = "17+38="
$t.Cell(19,5).Range.Text = "85-31="
$t.Cell(20,1).Range.Text = "81-74="
$t.Cell(20,2).Range.Text = "62+15="
$t.Cell(20,3).Range.Text = "76-12="
$t.Cell(20,4).Range.Text = "25+49="
$t.Cell(20,5).Range.Text = "21+63="
